# Insert a new "Task lần 3" section after the paragraph that contains
# "Thêm luật cho vấn đề sửa code người khác" (the last real content
# paragraph before the trailing blank paragraphs at the end of the doc).

$d = $word.ActiveDocument

# Locate the anchor paragraph by its text instead of a hard-coded index,
# so the script is resilient to any small structural differences.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Thêm luật cho vấn đề sửa code người khác*") {
        $anchor = $para
    }
}

$insPoint = $d.Range($anchor.Range.End, $anchor.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Task lần 3</w:t></w:r></w:p>' + `
'<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Cường</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Làm hệ thống Account và TopPlayer</w:t></w:r></w:p>' + `
'<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Tô Đức Anh</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Nghiên cứu J</w:t></w:r><w:r><w:t>U</w:t></w:r><w:r><w:t>nit</w:t></w:r></w:p>' + `
'<w:p/>' + `
'</w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$insPoint.InsertXML($xml)
